# Auto-generated edit script: update cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.231.24'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.997.11'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '544.54'
$ws.Range('E5').Value = '  -1.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.34'
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.993.40'
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('E9').Value = '  -1.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.89'
$ws.Range('E10').Value = '  +13.00%  '
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.445'
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.94'
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.471.78'
$ws.Range('E15').Value = '  -1.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.303.85'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.996.67'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.108'
$ws.Range('E18').Value = '  -2.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.55'
$ws.Range('E19').Value = '  -1.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '469.02'
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.37'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.652'
$ws.Range('E22').Value = '  -3.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.17'
$ws.Range('E23').Value = '  +1.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.30'
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.59'
$ws.Range('E25').Value = '  +3.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.64'
$ws.Range('E28').Value = '  -1.99%  '
$ws.Range('E29').Value = '  +4.34%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.37'
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.12'
$ws.Range('E32').Value = '  -2.90%  '
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '54.44'
$ws.Range('E35').Value = '  -2.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.83'
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '452.07'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0808'
$ws.Range('E38').Value = '  +1.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0391'
$ws.Range('E39').Value = '  +1.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.942.13'
$ws.Range('E40').Value = '  -8.24%  '
$ws.Range('E41').Value = '  -3.48%  '
$ws.Range('E42').Value = '  -1.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.52'
$ws.Range('E43').Value = '  +2.19%  '
$ws.Range('E44').Value = '  +3.24%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.249'
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.99'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.109'
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '115.15'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0492'
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('E51').Value = '  -0.88%  '
